$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds 13 brand rows (rows 2-14). The update expands this
# to 23 brand rows (rows 2-24): 10 new brands are inserted at the top (after
# the header) and the previously existing brands shift down by 10 rows, with
# a handful of their "id" (column C) values also changing.

# First, extend the formatting (style) that column A carries (bold/centered/
# bordered, style index 1 in the original file) down into the newly used
# rows 15-24, by copying the format from the last existing data row (A14).
$ws.Range("A14").Copy() | Out-Null
$ws.Range("A15:A24").PasteSpecial(-4122) | Out-Null

# Now populate every data row (2-24) with its final Brand/id values.
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "OPPO"
$ws.Cells.Item(2, 3).Value = 13

$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "Samsung"
$ws.Cells.Item(3, 3).Value = 17

$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "Realme"
$ws.Cells.Item(4, 3).Value = 16

$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "Original"
$ws.Cells.Item(5, 3).Value = 14

$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = "TCL"
$ws.Cells.Item(6, 3).Value = 18

$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = "Xiaomi"
$ws.Cells.Item(7, 3).Value = 20

$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "POCO"
$ws.Cells.Item(8, 3).Value = 15

$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "VIVO"
$ws.Cells.Item(9, 3).Value = 19

$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "Google"
$ws.Cells.Item(10, 3).Value = 8

$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "мир гаджетов"
$ws.Cells.Item(11, 3).Value = 22

$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "Acer"
$ws.Cells.Item(12, 3).Value = 1

$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "lenovo"
$ws.Cells.Item(13, 3).Value = 21

$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "Apple"
$ws.Cells.Item(14, 3).Value = 2

$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = "Gigabyte"
$ws.Cells.Item(15, 3).Value = 7

$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "Asus"
$ws.Cells.Item(16, 3).Value = 3

$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "HP"
$ws.Cells.Item(17, 3).Value = 9

$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = "Infinix"
$ws.Cells.Item(18, 3).Value = 12

$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = "Azerty"
$ws.Cells.Item(19, 3).Value = 4

$ws.Cells.Item(20, 1).Value = 18
$ws.Cells.Item(20, 2).Value = "Dere"
$ws.Cells.Item(20, 3).Value = 6

$ws.Cells.Item(21, 1).Value = 19
$ws.Cells.Item(21, 2).Value = "DIGMA"
$ws.Cells.Item(21, 3).Value = 5

$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(22, 2).Value = ""
$ws.Cells.Item(22, 3).Value = 0

$ws.Cells.Item(23, 1).Value = 21
$ws.Cells.Item(23, 2).Value = "Huawei"
$ws.Cells.Item(23, 3).Value = 10

$ws.Cells.Item(24, 1).Value = 22
$ws.Cells.Item(24, 2).Value = "IPASON"
$ws.Cells.Item(24, 3).Value = 11
